$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.110.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.432.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.65"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.592"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.44"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.955.16"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.413.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.194.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000132"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.44%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.59"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "316.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.12"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.79"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.18%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.68"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.40"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0484"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.68"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.291"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.135.10"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.87"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.04%  "
